# Update the "Overview" sheet: roll the quarterly data window forward by one
# quarter. Drop "Q2 ending 1399/06", shift every quarter's figures one
# column to the left (E<-F, F<-G, ... M<-N), and populate the newly
# freed-up last column (N) with the new quarter "Q4 ending 1401/12".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the quarter row (row 8 and row 24), columns E:N
$headers = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 5 + $i   # column E = 5
    $ws.Cells.Item(8, $col).Value = $headers[$i]
    $ws.Cells.Item(24, $col).Value = $headers[$i]
}

# Data rows: shift quarterly figures one column left and append the new
# quarter's values in column N. Ordered pairs (row, values) so write order
# is deterministic.
$dataRows = @(
    @{ Row = 10; Values = @(-4018, -825, 59, 794, -466, 496, 200, 262, 350, 218) },
    @{ Row = 11; Values = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0) },
    @{ Row = 12; Values = @(86720, 37555, 118282, 152247, 176274, 114962, 162693, 204610, 182490, 110698) },
    @{ Row = 13; Values = @(2687, 3682, 2597, 2908, 1455, 9216, 2751, 2103, 8013, 3200) },
    @{ Row = 14; Values = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0) },
    @{ Row = 15; Values = @(24, 50, 34, 40, 37, 32, 18, 22, 39, 40) },
    @{ Row = 16; Values = @(0, 963, 462, 722, 612, 902, 757, 1442, 1090, 1055) },
    @{ Row = 17; Values = @(28330, 27295, 34378, 50165, 31114, 48022, 53159, 68365, 48507, 51972) },
    @{ Row = 18; Values = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0) },
    @{ Row = 19; Values = @(3475, 25230, 19289, 10768, 25117, 38554, 10574, 20011, 11903, 25316) },
    @{ Row = 20; Values = @(117218, 93950, 175101, 217644, 234143, 212184, 230152, 296815, 252392, 192499) },
    @{ Row = 26; Values = @(366, 364, 365, 375, 375, 373, 369, 363, 366, 315) },
    @{ Row = 27; Values = @(521, 508, 513, 500, 502, 500, 523, 500, 502, 500) }
)

foreach ($entry in $dataRows) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i   # column E = 5
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
